$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Carlos Salinas bewell / csalinas@somosbewell.cl (replaces "Carlos Salinas Real" row)
$ws.Range("A2").Value = "Carlos Salinas bewell"
$ws.Range("B2").Value = "csalinas@somosbewell.cl"

# Update row 3: Claudio Cáceres / ccaceres@somosbewell.cl (replaces "Eduardo Reyes bewell" row)
$ws.Range("A3").Value = "Claudio Cáceres"
$ws.Range("B3").Value = "ccaceres@somosbewell.cl"

# Update row 4: Katherinne Campos / katherinne.campos@somosbewell.cl (new entry)
$ws.Range("A4").Value = "Katherinne Campos"
$ws.Range("B4").Value = "katherinne.campos@somosbewell.cl"

# Point each existing mailto hyperlink at the matching new e-mail address.
# Iterate the live collection (instead of indexing via .Item()) so each
# hyperlink is updated in place - re-using the existing hyperlink entry and
# relationship id rather than appending a new one.
$newAddresses = @{
    '$B$2' = "mailto:csalinas@somosbewell.cl"
    '$B$3' = "mailto:ccaceres@somosbewell.cl"
    '$B$4' = "mailto:katherinne.campos@somosbewell.cl"
}
foreach ($h in $ws.Hyperlinks) {
    $h.Address = $newAddresses[$h.Range.Address()]
}

# Update selection to C5
$ws.Range("C5").Select()
